$d = $word.ActiveDocument

# --- p086v: merge "<id>" + "p086v_a1" + "</id>" runs into a single
#     "<id>p086v_1</id>" run (keeps the first run's formatting). ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("<id>p086v_a1</id>", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = "<id>p086v_1</id>"
    Write-Output "p086v: <id>p086v_a1</id> -> <id>p086v_1</id>"
} else {
    Write-Output "p086v: pattern not found"
}

# --- p087r: merge "<id>" + "p087r_a1" + "</id>" runs into a single
#     "<id>p087r_1</id>" run (keeps the first run's formatting). ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("<id>p087r_a1</id>", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Text = "<id>p087r_1</id>"
    Write-Output "p087r: <id>p087r_a1</id> -> <id>p087r_1</id>"
} else {
    Write-Output "p087r: pattern not found"
}
